$p = $ppt.ActivePresentation

# --- Slide 12: re-position the background picture and drop the leftover
#     formula-annotation shapes (object 4 .. object 9) that used to sit on
#     top of it. ---
$s12 = $p.Slides.Item(12)

$pic12 = $s12.Shapes.Item("object 3")
$pic12.Left = 2.756142   # 35003 EMU
$pic12.Top = 0           # 0 EMU

$s12.Shapes.Item("object 4").Delete()
$s12.Shapes.Item("object 5").Delete()
$s12.Shapes.Item("object 6").Delete()
$s12.Shapes.Item("object 7").Delete()
$s12.Shapes.Item("object 8").Delete()
$s12.Shapes.Item("object 9").Delete()

# --- Slide 30: drop the same kind of leftover formula-annotation shapes. ---
$s30 = $p.Slides.Item(30)

$s30.Shapes.Item("object 5").Delete()
$s30.Shapes.Item("object 7").Delete()
$s30.Shapes.Item("object 8").Delete()
$s30.Shapes.Item("object 9").Delete()
$s30.Shapes.Item("object 10").Delete()
